$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Cells.Item(10,1).Value = @'
I/UX Designer - Figma Expert for Travel Platform 
Posted 1 hour ago 
Worldwide 
Specialized profiles can help you better highlight your expertise when submitting proposals to jobs like these. Create a specialized profile. 
Summary 
We're looking for a talented designer to elevate the user experience of our new website, including our cutting-edge AI travel itinerary planner and various platform interfaces. This role 
offers the opportunity to work with a mission-driven nonprofit while creating meaningful digital experiences for travelers worldwide. 
Objectives: 
● Design intuitive UI/UX for our AI-powered travel itinerary planning tool 
● Create compelling landing pages that convert visitors into engaged users 
● Develop seamless in-session page designs for optimal user journey flow 
● Design responsive interfaces across desktop, tablet, and mobile platforms 
● Support ad-hoc design projects as they arise within our growing organization 
Requirements: 
● Expert-level Figma proficiency - Advanced knowledge of Auto Layout, Variants, Interactive Components, and Variables 
● WordPress design experience - Understanding of WordPress themes, Elementor page builder, and plugin integration 
● Experience with Ninja Forms or similar form builders for seamless design integration 
● Proven experience designing responsive web applications and mobile interfaces 
● Strong understanding of user-centered design principles and conversion optimization 
● Portfolio showcasing complex responsive layouts and design system work 
Preferred: 
● Experience designing AI-powered interfaces or conversational UI patterns 
● Familiarity with travel planning workflows and user behavior 
● Background in cultural or religious tourism markets 
www.worldjewishtravel.org 
UNSOLICITED EMAILS OUTSIDE OF UPWORK WILL NOT BE CONSIDERED FOR THIS ROLE. 
Less than 30 hrs/week 
Hourly 
3 to 6 months 
Duration 
Expert 
I am willing to pay higher rates for the most experienced freelancers 
$25.00 
- 
$40.00 
Hourly 
Project Type: Ongoing project 
About the client 
Payment method verified 
Rating is 5.0 out of 5. 
5.0 
4.99 of 35 reviews 
United States 
Worcester 6:02 PM 
43 jobs posted 
96% hire rate, 3 open jobs 
$81K total spent 
49 hires, 3 active 
$27.01 /hr avg hourly rate paid 
2,551 hours 
Travel & Hospitality 
Small company (2-9 people) 
Member since Mar 28, 2018 
Job link
'@
$ws.Cells.Item(10,2).Value = @'
**Justification for Tone:**  
Given the client’s verified payment method, excellent 5-star rating, substantial budget, and the strategic importance of the project involving advanced Figma and travel industry expertise, a confident expert tone is ideal to convey solid experience, reliability, and professionalism.
---
**Upwork Proposal:**
Hello,
I’m excited about your search for a skilled UI/UX Designer and Figma expert to elevate your AI-powered travel itinerary platform. With a strong background in crafting intuitive and responsive travel and tourism interfaces, I understand the nuances that drive traveler engagement and conversion optimization.
To demonstrate my expertise, here are some relevant case studies and portfolio examples in travel platform UI/UX and Figma design:
- Upwork Case Study (UX/UI Design & Figma):  
  https://www.upwork.com/freelancers/muhammadadnan144?s=1044578476142100494&p=1709584746479128576
- Travel & Tourism Portfolio Projects:  
  https://insitechstaging.com/demo/travel-agency  
  https://insitechstaging.com/demo/relocate-travel
- Digital Agency Portfolio Showcasing Complex Responsive Design:  
  https://portfolio.muhammadprojects.com/voyage-v2/
Additionally, here is a Figma prototype showcasing advanced Figma features like Auto Layout, Interactive Components, Variants, and Variables within a travel UI context:  
[Figma Prototype - Redeemo](https://www.figma.com/file/4naoKxypi3BeeOQxWiQoot/Redeemo?type=design&node-id=0-1&mode=design)
**This is the website we have designed and the prototype link is above. This is what we have developed:**  
Using these advanced tools allows me to deliver seamless and easily maintainable design systems that fully support responsive, multi-device user journeys.
I am also well-versed in WordPress theme design, Elementor page builder, and Ninja Forms integration, ensuring your design vision smoothly translates into a robust platform.
**Project Timeline:**  
Typically, initial wireframes and landing page concepts can be delivered within 2 weeks, followed by interactive prototype iterations, with full responsive UI/UX designs for desktop, tablet, and mobile ready within 4-6 weeks depending on scope.
I’m eager to contribute to your mission-driven nonprofit by crafting a travel platform that not only looks compelling but drives engagement through user-centered design and conversion optimization.
Let’s connect to discuss your goals and how I can help bring this project to life!
Best regards,  
Muhammad Ibrahim
'@
$ws.Cells.Item(10,3).Value = @'
https://www.upwork.com/freelancers/muhammadadnan144?s=1044578476142100494&p=1709584746479128576
'@
$ws.Cells.Item(10,4).Value = @'
https://insitechstaging.com/demo/travel-agency
'@
$ws.Cells.Item(10,5).Value = @'
https://insitechstaging.com/demo/relocate-travel
'@
$ws.Cells.Item(10,6).Value = @'
https://portfolio.muhammadprojects.com/voyage-v2/
'@
$ws.Cells.Item(10,7).Value = @'
https://www.figma.com/file/4naoKxypi3BeeOQxWiQoot/Redeemo?type=design&node-id=0-1&mode=design
'@
$ws.Cells.Item(10,8).Value = 23951

# Row 11
$ws.Cells.Item(11,1).Value = @'
WordPress Expert Needed for Book Publication Landing Page
Posted 9 hours ago
Worldwide
Summary
We are seeking a skilled WordPress expert to design a visually appealing and user-friendly landing page for our upcoming book publication. We’re launching a new book titled "Reimagining Fairness: An Equity, Cultural Diversity, and Inclusion Competency Approach" and are seeking a skilled WordPress expert to create a modern, responsive landing page that effectively promotes the book and captures leads.
The ideal candidate will have experience in designing and developing landing pages that convert visitors into readers or subscribers. You should be proficient in customizing themes, optimizing for SEO, and ensuring mobile responsiveness. If you have a passion for books and a knack for creating engaging online experiences, we would love to hear from you!
We already have a Bluehost account with WordPress installed, and the project is ready to begin immediately.
Project Scope:
Design and develop a one-page WordPress landing page that includes:
* Book cover and visual branding - we provide
* Author bio and headshot - we provide
* Key takeaways or preview sections of the book - we provide
* Email sign-up integration (Jotform preferred - we provide)
* “Pre-order” and/or “Buy Now” buttons linked to Amazon and other retailers
* Testimonials or review quotes - we provide
* Social media links and share options - we provide
* SEO-friendly structure
* Mobile responsiveness and fast load times
* Built using Elementor (preferred), Divi, or native Gutenberg blocks
Requirements:
* Proven experience building high-converting landing pages for authors, products, or online courses
* Strong design sense and ability to work within branding guidelines
* Ability to recommend layout, structure, and conversion strategies
* Experience working with Bluehost-hosted WordPress sites
* Clear communication and timely delivery
Deliverables:
* Fully functional landing page live on our Bluehost-hosted WordPress site
* Instructions or brief Loom video for future edits
* Any stock images, plugins, or tools used (license compliant)
To Apply, Please Include:
* Links to 2–3 landing pages you've built (especially book or product-focused)
* Your preferred WordPress page builder and why
* Your estimated cost and turnaround time
* Any suggestions to make the page more effective or dynamic
Bonus Points For:
* Familiarity with DEI, social impact, or thought leadership projects
* Copywriting or layout suggestions that support lead generation and engagement
* Experience embedding video clips or animations
* We’re excited to collaborate with someone creative and reliable to help bring Reimagining Fairness to the world.
Less than 30 hrs/week
Hourly
Less than 1 month
Duration
Intermediate
I am looking for a mix of experience and value
$25.00
-
$45.00
Hourly
Project Type: One-time project
You will be asked to answer the following questions when submitting a proposal:
Describe your recent experience with similar projects
Please list any certifications related to this project
What frameworks have you worked with?
Include a link to your GitHub profile and/or website
Describe your approach to testing and improving QA
About the client 
Payment method verified 
Rating is 4.6 out of 5. 
4.6 
4.60 of 11 reviews 
United States 
San Francisco 12:05 PM 
34 jobs posted 
56% hire rate, 1 open job 
$138K total spent 
23 hires, 6 active 
$30.91 /hr avg hourly rate paid 
4,375 hours 
Education 
Small company (2-9 people)
'@
$ws.Cells.Item(11,2).Value = @'
Agent stopped due to iteration limit or time limit.
'@
$ws.Cells.Item(11,3).Value = ""
$ws.Cells.Item(11,4).Value = ""
$ws.Cells.Item(11,5).Value = ""
$ws.Cells.Item(11,6).Value = ""
$ws.Cells.Item(11,7).Value = ""
$ws.Cells.Item(11,8).Value = 26054

# Row 12
$ws.Cells.Item(12,1).Value = @'
Website Developer Needed to Fix Website Issues 
Posted 3 hours ago 
Only freelancers located in the U.S. may apply. 
Summary 
We are seeking an experienced developer to address and fix various issues on our website. The ideal candidate should have a strong background in troubleshooting, debugging, and improving website functionality. Tasks may include but are not limited to resolving errors, enhancing site speed, and ensuring mobile responsiveness. If you are detail-oriented and have a passion for creating seamless user experiences, we would love to hear from you! 
Hours to be determined 
Hourly 
1 to 3 months 
Duration 
Intermediate 
I am looking for a mix of experience and value 
$15.00 
- 
$30.00 
Hourly 
Project Type: Ongoing project 
You will be asked to answer the following questions when submitting a proposal: 
Describe your recent experience with similar projects 
Please list any certifications related to this project 
About the client 
Payment method verified 
Rating is 5.0 out of 5. 
5.0 
5.00 of 7 reviews 
United States 
Woodland Hills 12:06 PM 
14 jobs posted 
79% hire rate, 1 open job 
$7.2K total spent 
12 hires, 4 active 
$33.48 /hr avg hourly rate paid 
204 hours 
Mid-sized company (10-99 people) 

'@
$ws.Cells.Item(12,2).Value = @'
**Justification for tone:**  
Given the client’s preference for a mix of experience and value, combined with their mid-sized company stature paying average hourly rates above the job’s budget and 5-star rating, a confident expert tone is appropriate to demonstrate professionalism, technical competence, and reliability.
---
**Proposal:**  
Hello,
I’m excited to submit my proposal for your Website Developer role focused on fixing site issues, enhancing speed, and improving mobile responsiveness. With solid experience troubleshooting and debugging websites to ensure seamless user experiences, I’m confident I can provide the value and intermediate expertise you’re seeking.
Here are some examples from my Upwork portfolio that demonstrate my proficiency in website development and problem-solving:  
- [Construction B2B Website, UX/UI, & Service Architecture](https://www.upwork.com/freelancers/muhammadadnan144?s=1044578476142100494&p=1666125191209713664)
Additionally, here are some relevant portfolio pieces showing a range of website troubleshooting and functionality improvements:  
- [Forum33 - General Website Development](https://forum33.surge.sh/)  
- [Monsoon Security - Technology & IT](https://insitechstaging.com/demo/monsoon-security/)  
- [MyLres - Real Estate Property Management](https://mylres.com/)
To further underline my commitment to detailed and user-focused design, here is a Figma prototype showcasing a well-crafted responsive ecommerce website design:  
- [Vault Of Beauty Ecommerce Website Prototype](https://www.figma.com/design/2Khg5PQecG5oO6wWT7W45O/Vault-Of-Beauty?node-id=0-1&t=sVTTHwrxvuY2ZgQ0-1)
Though the job does not mention specific technologies or tools, I am proficient with all major web development and debugging technologies and tools required to identify and resolve errors effectively, optimize site speed, and ensure mobile responsiveness.
Based on your job description, I estimate it would take about 1-3 months to thoroughly diagnose, fix issues, and enhance your website for optimal performance, matching your indicated timeline and hourly budget.
I look forward to the opportunity to discuss your project in detail and demonstrate how I can help deliver seamless website functionality. Let’s connect and take the next step!
Best regards,  
Muhammad Ibrahim
'@
$ws.Cells.Item(12,3).Value = @'
https://www.upwork.com/freelancers/muhammadadnan144?s=1044578476142100494&p=1666125191209713664
'@
$ws.Cells.Item(12,4).Value = @'
https://forum33.surge.sh/
'@
$ws.Cells.Item(12,5).Value = @'
https://insitechstaging.com/demo/monsoon-security/
'@
$ws.Cells.Item(12,6).Value = @'
https://mylres.com/
'@
$ws.Cells.Item(12,7).Value = @'
https://www.figma.com/design/2Khg5PQecG5oO6wWT7W45O/Vault-Of-Beauty?node-id=0-1&t=sVTTHwrxvuY2ZgQ0-1
'@
$ws.Cells.Item(12,8).Value = 16681

# Row 13
$ws.Cells.Item(13,1).Value = @'
Web Designer needed for simple 4 page website
Posted 3 hours ago
Worldwide
Summary
We're looking for a creative website designer to work with us on created a new 4 page website for one of our new products.
The structure of the pages will be provided by us along logo, colours and font + all copy. The pages needed are...
Homepage
About Us
Pricing
Contact Form
Please share examples of your previous website design work.
Only apply via UpWork, any contact or proposal outside of Upwork will be ignored.
Web programming will be handed by our developer, design only required in Figma or similar.
Thanks!
Hours to be determined
Hourly
1 to 3 months
Duration
Intermediate
I am looking for a mix of experience and value
$15.00
-
$30.00
Hourly
Project Type: Ongoing project
You will be asked to answer the following questions when submitting a proposal:
Include examples of recent website design work
About the client
Payment method verified
Rating is 4.9 out of 5.
4.9
4.89 of 145 reviews
United Kingdom

'@
$ws.Cells.Item(13,2).Value = @'
Agent stopped due to iteration limit or time limit.
'@
$ws.Cells.Item(13,3).Value = ""
$ws.Cells.Item(13,4).Value = ""
$ws.Cells.Item(13,5).Value = ""
$ws.Cells.Item(13,6).Value = ""
$ws.Cells.Item(13,7).Value = ""
$ws.Cells.Item(13,8).Value = 16518

# Row 14
$ws.Cells.Item(14,1).Value = @'
Landing Page Specialist Needed for Funnel 
Posted 3 hours ago 
Worldwide 
Summary 
We are a high-performance Marketing Agency working with fast-scaling E-Commerce and D2C brands. We’re looking for a top-tier designer who can translate our bold, results-driven brand into sleek visual assets – with a focus on clarity, professionalism, and conversion. 
What we need: 
1. 3 Conversion-Optimized Landing Page Designs (Funnel) 
2. Email Signature (responsive & clean) 
3. Business Cards (double-sided, premium design) 
– Each tailored to a specific buyer persona 
– Structure-first, mobile-first, built for lead generation 
– No cookie-cutter designs – we need premium, high-ticket visuals 
Target Audience: 
E-Com brands (6–7 figure monthly revenue), D2C founders, Heads of Marketing. 
Visual Style: Bold. Clear. Premium. No playful startup aesthetics. 
What we’re NOT looking for: 
– No generic templates 
– No stock-heavy visuals 
– No design without understanding conversion logic 
What we expect from you: 
– Experience in branding & conversion-focused design 
– Strong understanding of performance marketing funnels 
– Fast turnaround & clear communication 
– Delivery of all source files (Figma, PDF, PNG, SVG, etc.) 
Bonus if you have: 
– Webflow experience or dev support in your team 
– Prior work with D2C / E-Commerce brands 
How to Apply: 
Please send us: 
✅ 3–5 relevant portfolio pieces (please no generic portfolios) 
✅ Your price estimate for the full project 
✅ Timeline & workflow 
✅ Short message on why you’re the perfect fit 
Let’s build a bold, high-converting brand together. 
Looking forward to your application! 🚀 
More than 30 hrs/week 
Hourly 
1 to 3 months 
Duration 
Intermediate 
I am looking for a mix of experience and value 
$15.00 
- 
$40.00 
Hourly 
Contract-to-hire opportunity 
This lets talent know that this job could become full time. 
Learn more 
Project Type: Ongoing projectAbout the client 
Payment method verified 
Phone number verified 
Rating is 5.0 out of 5. 
5.0 
5.00 of 7 reviews 
United States 
Cleveland 11:09 AM 
42 jobs posted 
41% hire rate, 2 open jobs 
$68K total spent 
18 hires, 11 active 
$39.02 /hr avg hourly rate paid 
1,543 hours 
Small company (2-9 people)
'@
$ws.Cells.Item(14,2).Value = @'
Agent stopped due to iteration limit or time limit.
'@
$ws.Cells.Item(14,3).Value = ""
$ws.Cells.Item(14,4).Value = ""
$ws.Cells.Item(14,5).Value = ""
$ws.Cells.Item(14,6).Value = ""
$ws.Cells.Item(14,7).Value = ""
$ws.Cells.Item(14,8).Value = 25796
